$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new price value is a plain decimal number: force Text format
# first so Excel keeps the original text representation (matching trailing
# zeros, etc.) instead of auto-converting the cell to a numeric value.
$textCells = @("D5", "D6", "D9", "D15", "D16", "D19", "D20", "D23", "D25", "D29", "D37", "D46", "D47", "D48", "D50")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '26.190.94'
$ws.Range("E2").Value = '  -0.45%  '
$ws.Range("D3").Value = '1.589.95'
$ws.Range("E3").Value = '  +0.07%  '
$ws.Range("E4").Value = '  -0.12%  '
$ws.Range("D5").Value = '211.77'
$ws.Range("E5").Value = '  +0.81%  '
$ws.Range("D6").Value = '0.502'
$ws.Range("E6").Value = '  -0.91%  '
$ws.Range("E7").Value = '  -0.13%  '
$ws.Range("E8").Value = '  -0.16%  '
$ws.Range("D9").Value = '0.0605'
$ws.Range("E9").Value = '  -0.89%  '
$ws.Range("E10").Value = '  -2.12%  '
$ws.Range("E11").Value = '  -0.11%  '
$ws.Range("D12").Value = '1.812.63'
$ws.Range("E12").Value = '  +0.01%  '
$ws.Range("D13").Value = '1.582.85'
$ws.Range("E13").Value = '  -0.21%  '
$ws.Range("E14").Value = '  -1.39%  '
$ws.Range("D15").Value = '0.510'
$ws.Range("E15").Value = '  -1.62%  '
$ws.Range("D16").Value = '63.63'
$ws.Range("E16").Value = '  -1.11%  '
$ws.Range("D17").Value = '26.167.08'
$ws.Range("E17").Value = '  -0.57%  '
$ws.Range("D18").Value = '0.0₃0724'
$ws.Range("E18").Value = '  -0.58%  '
$ws.Range("D19").Value = '214.97'
$ws.Range("E19").Value = '  +1.93%  '
$ws.Range("D20").Value = '7.36'
$ws.Range("E20").Value = '  -1.74%  '
$ws.Range("E21").Value = '  -0.15%  '
$ws.Range("E22").Value = '  -0.61%  '
$ws.Range("D23").Value = '9.01'
$ws.Range("E23").Value = '  +0.75%  '
$ws.Range("E24").Value = '  -1.98%  '
$ws.Range("D25").Value = '144.54'
$ws.Range("E25").Value = '  -0.01%  '
$ws.Range("E26").Value = '  -0.11%  '
$ws.Range("E27").Value = '  -1.15%  '
$ws.Range("E28").Value = '  -0.94%  '
$ws.Range("D29").Value = '15.08'
$ws.Range("E29").Value = '  -1.08%  '
$ws.Range("E30").Value = '  -2.37%  '
$ws.Range("E31").Value = '  +0.22%  '
$ws.Range("E32").Value = '  -1.53%  '
$ws.Range("D33").Value = '1.417.85'
$ws.Range("E33").Value = '  +7.99%  '
$ws.Range("E34").Value = '  -1.50%  '
$ws.Range("E35").Value = '  -0.75%  '
$ws.Range("E36").Value = '  -0.64%  '
$ws.Range("D37").Value = '0.585'
$ws.Range("E37").Value = '  -4.31%  '
$ws.Range("E38").Value = '  -1.80%  '
$ws.Range("E39").Value = '  +2.12%  '
$ws.Range("E40").Value = '  +4.52%  '
$ws.Range("E42").Value = '  -13.92%  '
$ws.Range("E43").Value = '  +0.39%  '
$ws.Range("E44").Value = '  -0.71%  '
$ws.Range("D45").Value = '1.724.83'
$ws.Range("E45").Value = '  -0.01%  '
$ws.Range("D46").Value = '61.09'
$ws.Range("E46").Value = '  -1.84%  '
$ws.Range("D47").Value = '87.37'
$ws.Range("D48").Value = '1.49'
$ws.Range("E48").Value = '  -0.36%  '
$ws.Range("E49").Value = '  -0.82%  '
$ws.Range("D50").Value = '0.0960'
$ws.Range("E50").Value = '  -2.14%  '
$ws.Range("E51").Value = '  -0.18%  '
